# NYPD CompStat weekly report update ("New crime data collected")
# Applies the numeric/text updates described by the source diff to the
# already-open workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich-text shared strings).
# Use Characters() so only the targeted substring is rewritten, leaving the
# rest of the surrounding text (and cell) untouched.
# ---------------------------------------------------------------------------

# A8 holds: "Volume 32   Number  29" -> bump the report number 29 -> 30
$volRange = $ws.Range("A8")
$volText = $volRange.Value()
$numPos = $volText.IndexOf("29") + 1
$volRange.Characters($numPos, 2).Text = "30"

# C9 holds: "Report Covering the Week  7/14/2025  Through  7/20/2025"
# -> shift both dates forward by one week
$weekRange = $ws.Range("C9")
$weekText = $weekRange.Value()
$pos1 = $weekText.IndexOf("7/14/2025") + 1
$weekRange.Characters($pos1, 9).Text = "7/21/2025"
$pos2 = $weekText.IndexOf("7/20/2025") + 1
$weekRange.Characters($pos2, 9).Text = "7/27/2025"

# ---------------------------------------------------------------------------
# Row 15 (Rape)
# ---------------------------------------------------------------------------
$ws.Range("F15").Value = 1

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = -71.428571428571
$ws.Range("J16").Value = 37
$ws.Range("K16").Value = -27.027027027027
$ws.Range("M16").Value = -27.027027027027
$ws.Range("N16").Value = -82.692307692307

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = 59
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = -41.584158415841
$ws.Range("L17").Value = -9.230769230769
$ws.Range("M17").Value = 11.320754716981
$ws.Range("N17").Value = -51.639344262295

# ---------------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 36.363636363636
$ws.Range("M18").Value = -11.764705882352
$ws.Range("N18").Value = -88.593155893536

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -71.428571428571
$ws.Range("F19").Value = 13
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = -27.777777777777
$ws.Range("I19").Value = 72
$ws.Range("J19").Value = 90
$ws.Range("K19").Value = -20
$ws.Range("L19").Value = -34.545454545454
$ws.Range("M19").Value = 14.285714285714
$ws.Range("N19").Value = -38.983050847457

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------------
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 15
$ws.Range("J20").Value = 32
$ws.Range("K20").Value = -53.125
$ws.Range("L20").Value = -31.818181818181
$ws.Range("M20").Value = 7.142857142857
$ws.Range("N20").Value = -92.537313432835

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -57.142857142857
$ws.Range("F21").Value = 35
$ws.Range("G21").Value = 46
$ws.Range("H21").Value = -23.913043478260
$ws.Range("I21").Value = 207
$ws.Range("J21").Value = 293
$ws.Range("K21").Value = -29.351535836177
$ws.Range("L21").Value = -14.462809917355
$ws.Range("M21").Value = 0.975609756097
$ws.Range("N21").Value = -76.234213547646

# ---------------------------------------------------------------------------
# Row 22 (Transit) - C22/D22/E22 become the "no data" placeholders
# ("0" / "***.*"), matching the style+text already used elsewhere in the
# sheet (e.g. row 23). Copy from an identical existing cell so the cell
# keeps the exact shared-string + style pairing instead of Excel
# re-interpreting the literal text as a number.
# ---------------------------------------------------------------------------
$ws.Range("C23").Copy($ws.Range("C22"))
$ws.Range("D23").Copy($ws.Range("D22"))
$ws.Range("E23").Copy($ws.Range("E22"))
$ws.Range("M22").Value = -76.923076923076

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 43
$ws.Range("H24").Value = 16.279069767441
$ws.Range("I24").Value = 292
$ws.Range("J24").Value = 253
$ws.Range("K24").Value = 15.415019762845
$ws.Range("L24").Value = 7.749077490774
$ws.Range("M24").Value = 81.366459627329

# ---------------------------------------------------------------------------
# Row 25 (Retail Theft)
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 19.047619047619
$ws.Range("I25").Value = 148
$ws.Range("J25").Value = 123
$ws.Range("K25").Value = 20.325203252032
$ws.Range("L25").Value = 1.369863013698

# ---------------------------------------------------------------------------
# Row 26 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -11.111111111111
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = -34.285714285714
$ws.Range("I26").Value = 131
$ws.Range("J26").Value = 158
$ws.Range("K26").Value = -17.088607594936
$ws.Range("L26").Value = -6.428571428571
$ws.Range("M26").Value = -36.097560975609

# ---------------------------------------------------------------------------
# Row 27 (UCR Rape*) - C27 becomes the "no data" placeholder ("0")
# ---------------------------------------------------------------------------
$ws.Range("C23").Copy($ws.Range("C27"))
$ws.Range("F27").Value = 2

# ---------------------------------------------------------------------------
# Row 28 (Other Sex Crimes) - C28/D28/E28 become "no data" placeholders
# ---------------------------------------------------------------------------
$ws.Range("C23").Copy($ws.Range("C28"))
$ws.Range("D23").Copy($ws.Range("D28"))
$ws.Range("E23").Copy($ws.Range("E28"))
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 300

# ---------------------------------------------------------------------------
# Row 29 (Shooting Vic.) - F29 becomes the "no data" placeholder ("0")
# ---------------------------------------------------------------------------
$ws.Range("F14").Copy($ws.Range("F29"))

# ---------------------------------------------------------------------------
# Row 30 (Shooting Inc.) - F30 becomes the "no data" placeholder ("0")
# ---------------------------------------------------------------------------
$ws.Range("F14").Copy($ws.Range("F30"))
